# Apply the "LinuxForHealth" rebrand + regen edit to the provider-region
# StructureDefinition workbook.
#
# Changes:
#  1. Metadata sheet: URL, Version, Date, Publisher values are refreshed to
#     reflect the move from Alvearie/ibm.com -> LinuxForHealth.
#  2. Elements sheet: the top-level "Extension" row's Constraint(s) cell
#     (which had erroneously carried the ele-1/ext-1 FHIRPath constraint
#     text meant for the "Extension.extension" row) is cleared out, as
#     happens when the IG publisher regenerates this view.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(2, 2).Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/provider-region"
$meta.Cells.Item(3, 2).Value = "8.0.0"
$meta.Cells.Item(8, 2).Value = "2022-11-10T16:00:46+00:00"
$meta.Cells.Item(9, 2).Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the "Extension" element; column AI (35) is "Constraint(s)".
$elements.Cells.Item(2, 35).Value = ""
